$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume cells are stored as plain text (e.g. "56.540.61",
# European-style thousands separators, or percentages with padding spaces).
# A naive `.Value = "498.25"` assignment lets Excel auto-coerce single-dot
# numeric-looking strings into real numbers, which would change the cell's
# stored type. Detect that coercion and force the text back in as a string
# while stripping the temporary Text number-format so the cell style is left
# exactly as it was (no lingering "@" format / style index).
function Set-TextValue($cell, $text) {
    $cell.Value = $text
    if ($cell.Value -ne $text) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    }
}

Set-TextValue $ws.Cells.Item(2, 4) '56.540.61'
Set-TextValue $ws.Cells.Item(2, 5) '  -2.87%  '
Set-TextValue $ws.Cells.Item(3, 4) '2.991.06'
Set-TextValue $ws.Cells.Item(3, 5) '  -4.75%  '
Set-TextValue $ws.Cells.Item(5, 4) '498.25'
Set-TextValue $ws.Cells.Item(5, 5) '  -5.10%  '
Set-TextValue $ws.Cells.Item(6, 4) '134.91'
Set-TextValue $ws.Cells.Item(6, 5) '  -0.06%  '
Set-TextValue $ws.Cells.Item(7, 5) '  -0.08%  '
Set-TextValue $ws.Cells.Item(8, 4) '2.988.12'
Set-TextValue $ws.Cells.Item(8, 5) '  -4.75%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.426'
Set-TextValue $ws.Cells.Item(9, 5) '  -4.09%  '
Set-TextValue $ws.Cells.Item(10, 4) '7.25'
Set-TextValue $ws.Cells.Item(10, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.105'
Set-TextValue $ws.Cells.Item(11, 5) '  -3.55%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.352'
Set-TextValue $ws.Cells.Item(12, 5) '  -7.62%  '
Set-TextValue $ws.Cells.Item(14, 4) '3.496.76'
Set-TextValue $ws.Cells.Item(14, 5) '  -4.84%  '
Set-TextValue $ws.Cells.Item(15, 4) '24.91'
Set-TextValue $ws.Cells.Item(15, 5) '  -2.56%  '
Set-TextValue $ws.Cells.Item(16, 4) '56.549.64'
Set-TextValue $ws.Cells.Item(16, 5) '  -2.78%  '
Set-TextValue $ws.Cells.Item(17, 5) '  -3.52%  '
Set-TextValue $ws.Cells.Item(18, 4) '2.987.10'
Set-TextValue $ws.Cells.Item(18, 5) '  -4.80%  '
Set-TextValue $ws.Cells.Item(19, 4) '5.84'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.67%  '
Set-TextValue $ws.Cells.Item(20, 4) '12.32'
Set-TextValue $ws.Cells.Item(20, 5) '  -5.95%  '
Set-TextValue $ws.Cells.Item(21, 4) '7.77'
Set-TextValue $ws.Cells.Item(21, 5) '  -2.45%  '
Set-TextValue $ws.Cells.Item(22, 4) '327.08'
Set-TextValue $ws.Cells.Item(22, 5) '  -5.17%  '
Set-TextValue $ws.Cells.Item(23, 5) '  +0.07%  '
Set-TextValue $ws.Cells.Item(24, 4) '0.466'
Set-TextValue $ws.Cells.Item(24, 5) '  -8.37%  '
Set-TextValue $ws.Cells.Item(25, 4) '61.31'
Set-TextValue $ws.Cells.Item(25, 5) '  -10.47%  '
Set-TextValue $ws.Cells.Item(26, 4) '0.999'
Set-TextValue $ws.Cells.Item(26, 5) '  -0.15%  '
Set-TextValue $ws.Cells.Item(27, 4) '0.161'
Set-TextValue $ws.Cells.Item(27, 5) '  -3.88%  '
Set-TextValue $ws.Cells.Item(28, 4) '0.0₃0913'
Set-TextValue $ws.Cells.Item(28, 5) '  -4.77%  '
Set-TextValue $ws.Cells.Item(29, 5) '  -0.06%  '
Set-TextValue $ws.Cells.Item(30, 5) '  -4.60%  '
Set-TextValue $ws.Cells.Item(31, 4) '6.87'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.06%  '
Set-TextValue $ws.Cells.Item(32, 4) '1.18'
Set-TextValue $ws.Cells.Item(32, 5) '  -3.74%  '
Set-TextValue $ws.Cells.Item(33, 4) '1.74'
Set-TextValue $ws.Cells.Item(33, 5) '  -6.78%  '
Set-TextValue $ws.Cells.Item(34, 4) '20.02'
Set-TextValue $ws.Cells.Item(34, 5) '  -7.03%  '
Set-TextValue $ws.Cells.Item(35, 4) '154.17'
Set-TextValue $ws.Cells.Item(35, 5) '  -1.86%  '
Set-TextValue $ws.Cells.Item(36, 4) '4.49'
Set-TextValue $ws.Cells.Item(36, 5) '  -6.82%  '
Set-TextValue $ws.Cells.Item(37, 5) '  -7.05%  '
Set-TextValue $ws.Cells.Item(38, 4) '5.61'
Set-TextValue $ws.Cells.Item(38, 5) '  -10.00%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.0677'
Set-TextValue $ws.Cells.Item(39, 5) '  -2.15%  '
Set-TextValue $ws.Cells.Item(40, 4) '23.61'
Set-TextValue $ws.Cells.Item(40, 5) '  -3.82%  '
Set-TextValue $ws.Cells.Item(41, 4) '3.021.06'
Set-TextValue $ws.Cells.Item(41, 5) '  -4.69%  '
Set-TextValue $ws.Cells.Item(42, 4) '36.67'
Set-TextValue $ws.Cells.Item(42, 5) '  -9.11%  '
Set-TextValue $ws.Cells.Item(43, 5) '  +0.08%  '
Set-TextValue $ws.Cells.Item(44, 5) '  -7.04%  '
Set-TextValue $ws.Cells.Item(45, 4) '0.638'
Set-TextValue $ws.Cells.Item(45, 5) '  -7.78%  '
Set-TextValue $ws.Cells.Item(46, 5) '  -2.17%  '
Set-TextValue $ws.Cells.Item(47, 4) '2.205.88'
Set-TextValue $ws.Cells.Item(47, 5) '  -2.53%  '
Set-TextValue $ws.Cells.Item(48, 4) '3.57'
Set-TextValue $ws.Cells.Item(48, 5) '  -8.72%  '
Set-TextValue $ws.Cells.Item(49, 5) '  +6.05%  '
Set-TextValue $ws.Cells.Item(50, 5) '  +1.52%  '
Set-TextValue $ws.Cells.Item(51, 4) '5.72'
Set-TextValue $ws.Cells.Item(51, 5) '  -7.64%  '
